# Insert two new rows at position 999 (shifting all following rows down by 2)
# and populate them with the new "Cebollín" price entries for fecha 44746
# (Primera and Segunda quality).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(999).Insert()
$ws.Rows.Item(999).Insert()

# New row 999: Primera
$ws.Cells.Item(999, 1).Value = 9
$ws.Cells.Item(999, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(999, 3).Value = "Metropolitana"
$ws.Cells.Item(999, 4).Value = 44746
$ws.Cells.Item(999, 5).Value = 13
$ws.Cells.Item(999, 6).Value = 100112037
$ws.Cells.Item(999, 7).Value = "Cebollín"
$ws.Cells.Item(999, 8).Value = "Sin especificar"
$ws.Cells.Item(999, 9).Value = "Primera"
$ws.Cells.Item(999, 10).Value = 160
$ws.Cells.Item(999, 11).Value = 7000
$ws.Cells.Item(999, 12).Value = 7000
$ws.Cells.Item(999, 13).Value = 7000
$ws.Cells.Item(999, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(999, 15).Value = "Región Metropolitana"
$ws.Cells.Item(999, 16).Value = 194
$ws.Cells.Item(999, 17).Value = 36
$ws.Cells.Item(999, 18).Value = "Hortaliza"

# New row 1000: Segunda
$ws.Cells.Item(1000, 1).Value = 9
$ws.Cells.Item(1000, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1000, 3).Value = "Metropolitana"
$ws.Cells.Item(1000, 4).Value = 44746
$ws.Cells.Item(1000, 5).Value = 13
$ws.Cells.Item(1000, 6).Value = 100112037
$ws.Cells.Item(1000, 7).Value = "Cebollín"
$ws.Cells.Item(1000, 8).Value = "Sin especificar"
$ws.Cells.Item(1000, 9).Value = "Segunda"
$ws.Cells.Item(1000, 10).Value = 70
$ws.Cells.Item(1000, 11).Value = 6000
$ws.Cells.Item(1000, 12).Value = 6000
$ws.Cells.Item(1000, 13).Value = 6000
$ws.Cells.Item(1000, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(1000, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1000, 16).Value = 167
$ws.Cells.Item(1000, 17).Value = 36
$ws.Cells.Item(1000, 18).Value = "Hortaliza"
